$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header value changes
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data value changes
$ws.Range("B2").Value = 50.720884012982765
$ws.Range("C2").Value = 56.772752221374418
$ws.Range("D2").Value = 53.094941844320054
$ws.Range("E2").Value = 60.124740101386323

# Row 3 data value changes
$ws.Range("B3").Value = 45.807443843297833
$ws.Range("C3").Value = 51.026369555577283
$ws.Range("D3").Value = 51.765353212861513
$ws.Range("E3").Value = 56.815140331805843

# Update selection to match new sqref B1:E3
$ws.Range("B1:E3").Select()
